$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Mã học sinh" column (column A), shifting everything left.
$ws.Columns.Item(1).Delete()

# 2. Fix up data values that changed in the new template.
#    E (Khối) is now a simple numeric 1 / 5 instead of 6 / 5.
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 5

#    F (Lớp) values changed from "6A1"/"5B2" to "1A"/"5B".
$ws.Range("F2").Value = "1A"
$ws.Range("F3").Value = "5B"

# 3. Student code (HS001 / HS002) is no longer a column of its own; the
#    student email cell (column B, row 2) now carries the student's own
#    mailto hyperlink as well as the existing parent-email hyperlinks
#    (which shifted from column J to column I after the column delete).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:HS001@email.com") | Out-Null

$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:parentA@email.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:parentB@email.com") | Out-Null
$ws.Range("I2").Interior.Pattern = -4142
$ws.Range("I3").Interior.Pattern = -4142

# 4. H (SĐT phụ huynh) now stores the phone numbers as zero-padded text, and the
#    whole column (header included) switches to a text number format.
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H2").Value = "0912345678"
$ws.Range("H3").Value = "0987654321"

Write-Host "base edits done"
